# Auto-generated script applying cryptos.xlsx data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.067.66"
$ws.Range("E2").Value = "  +2.19%  "
$ws.Range("D3").Value = "2.431.11"
$ws.Range("E3").Value = "  +1.97%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'307.94"
$ws.Range("E5").Value = "  +1.53%  "
$ws.Range("D6").Value = "'98.35"
$ws.Range("E6").Value = "  +0.75%  "
$ws.Range("D7").Value = "'0.512"
$ws.Range("E7").Value = "  +0.36%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").Value = "'0.498"
$ws.Range("E9").Value = "  -0.81%  "
$ws.Range("D10").Value = "'35.05"
$ws.Range("E10").Value = "  +2.31%  "
$ws.Range("D11").Value = "'0.0801"
$ws.Range("E11").Value = "  +1.36%  "
$ws.Range("E12").Value = "  +2.74%  "
$ws.Range("E13").Value = "  +0.49%  "
$ws.Range("D14").Value = "'6.95"
$ws.Range("E14").Value = "  +2.24%  "
$ws.Range("D15").Value = "2.802.35"
$ws.Range("E15").Value = "  +1.67%  "
$ws.Range("D16").Value = "2.464.05"
$ws.Range("E16").Value = "  +2.65%  "
$ws.Range("D17").Value = "'0.831"
$ws.Range("E17").Value = "  +2.10%  "
$ws.Range("D18").Value = "44.006.98"
$ws.Range("E18").Value = "  +2.05%  "
$ws.Range("B19").Value = "InternetComputer(DFINITY)"
$ws.Range("C19").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D19").Value = "'12.24"
$ws.Range("E19").Value = "  -0.75%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").Value = "'6.46"
$ws.Range("E20").Value = "  +1.16%  "
$ws.Range("D21").Value = "0.0₃0907"
$ws.Range("E21").Value = "  +1.98%  "
$ws.Range("D22").Value = "'68.32"
$ws.Range("E22").Value = "  -0.14%  "
$ws.Range("D23").Value = "'239.41"
$ws.Range("E23").Value = "  +1.36%  "
$ws.Range("E24").Value = "  +2.44%  "
$ws.Range("D25").Value = "'2.46"
$ws.Range("E25").Value = "  +0.96%  "
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("D27").Value = "'25.16"
$ws.Range("E27").Value = "  +1.43%  "
$ws.Range("D28").Value = "'2.20"
$ws.Range("E28").Value = "  -6.76%  "
$ws.Range("D29").Value = "'9.47"
$ws.Range("E29").Value = "  +3.68%  "
$ws.Range("D30").Value = "'32.74"
$ws.Range("E30").Value = "  +3.34%  "
$ws.Range("D31").Value = "'0.118"
$ws.Range("E31").Value = "  +15.80%  "
$ws.Range("D32").Value = "'18.60"
$ws.Range("E32").Value = "  +7.96%  "
$ws.Range("E33").Value = "  +0.96%  "
$ws.Range("D34").Value = "'1.00"
$ws.Range("E34").Value = "  +0.05%  "
$ws.Range("E35").Value = "  +3.72%  "
$ws.Range("E36").Value = "  +3.30%  "
$ws.Range("D37").Value = "'130.57"
$ws.Range("E37").Value = "  +21.26%  "
$ws.Range("D38").Value = "'4.46"
$ws.Range("E38").Value = "  +1.76%  "
$ws.Range("D39").Value = "'2.92"
$ws.Range("E39").Value = "  +4.11%  "
$ws.Range("E40").Value = "  -1.07%  "
$ws.Range("E41").Value = "  +0.14%  "
$ws.Range("D42").Value = "'21.29"
$ws.Range("E42").Value = "  -4.45%  "
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").Value = "'0.0286"
$ws.Range("E43").Value = "  +2.03%  "
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "1.950.29"
$ws.Range("E44").Value = "  -0.10%  "
$ws.Range("E45").Value = "  +2.15%  "
$ws.Range("E46").Value = "  +4.08%  "
$ws.Range("D47").Value = "'9.34"
$ws.Range("E47").Value = "  +0.49%  "
$ws.Range("D48").Value = "2.661.63"
$ws.Range("E48").Value = "  +1.95%  "
$ws.Range("D49").Value = "'1.61"
$ws.Range("E49").Value = "  +6.44%  "
$ws.Range("D50").Value = "'53.25"
$ws.Range("E50").Value = "  +0.49%  "
$ws.Range("D51").Value = "'73.01"
$ws.Range("E51").Value = "  +1.05%  "
